$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The letterhead rows (1-3) are removed by shifting the data table (rows 4-16)
# up by three rows. Row heights are a row-level property and must stay put,
# so this is done as cell-range cut/paste (not a true row delete), performed
# in small chunks because the COM shim's Cut(destination) only reliably
# relocates a few rows at a time.
$ws.Range("A4:F6").Cut($ws.Range("A1"))
$ws.Range("A7:F9").Cut($ws.Range("A4"))
$ws.Range("A10:F12").Cut($ws.Range("A7"))
$ws.Range("A13:F15").Cut($ws.Range("A10"))
$ws.Range("A16:F16").Cut($ws.Range("A13"))

# The vacated rows (now 14-16, previously holding Rheinland-Pfalz / Saarland /
# Schleswig-Holstein before they moved up) are fully cleared - values and
# formatting - leaving plain blank rows.
$ws.Range("A14:F16").Clear()

# The old footnote row (17) stays put but its text is removed while its
# cell formatting is kept.
$ws.Range("A17:B17").ClearContents()
